{"js": "const body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(\n  \"Systems of Water Supply and Treatment\",\n  \"Systems of Water Treatment\"\n);\n\nawait replaceOnce(\n  \"Ativa\u00e7\u00e3o: 01/01/2024\",\n  \"Ativa\u00e7\u00e3o: 01/01/2025\"\n);\n\nawait replaceOnce(\n  \"Fornecer aos alunos os conhecimentos b\u00e1sicos dos sistemas de abastecimento de \u00e1gua e das tecnologias de tratamento de \u00e1gua para consumo humano.\",\n  \"Fornecer aos alunos os conhecimentos b\u00e1sicos dos sistemas principais de tratamento de \u00e1gua para consumo humano.\"\n);\n\nawait replaceOnce(\n  \"Supply the students the basic knowledge of the systems of water supply and the water treatment technologies for the human consumption.\",\n  \"Supply the students the basic knowledge of the systems of main water treatment technologies for the human consumption.\"\n);\n\nawait replaceOnce(\n  \"Sistemas de Abastecimento de \u00c1gua; Tecnologias de Tratamento de \u00c1gua;\",\n  \"Tecnologias de Tratamento de \u00c1gua;\"\n);\n\nawait replaceOnce(\n  \"Systems of Water Supply; Water Treatment Technologies;\",\n  \"Water Treatment Technologies;\"\n);\n\nawait replaceOnce(\n  \"- Tipos de dimensionamento de redes de distribui\u00e7\u00e3o de \u00e1gua;- Caracter\u00edsticas\",\n  \"- Caracter\u00edsticas\"\n);\n\nawait replaceOnce(\n  \"- Reserva\u00e7\u00e3o; - Redes de distribui\u00e7\u00e3o: tipos de rede, pe\u00e7as e \u00f3rg\u00e3os acess\u00f3rios; - Sistema de Tratamento de \u00c1gua de Ciclo Completo;\",\n  \"- Reserva\u00e7\u00e3o; - Sistema de Tratamento de \u00c1gua de Ciclo Completo;\"\n);\n\nawait replaceOnce(\n  \"- Types of dimensioning of water distribution networks;- Water characteristics\",\n  \"- Water characteristics\"\n);\n\nawait replaceOnce(\n  \"Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;\",\n  \"Reservation;- Full Cycle Water Treatment System;\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute(\n        $find,           # FindText\n        $true,           # MatchCase\n        $false,          # MatchWholeWord\n        $false,          # MatchWildcards\n        $false,          # MatchSoundsLike\n        $false,          # MatchAllWordForms\n        $true,           # Forward\n        1,               # Wrap (wdFindContinue)\n        $false,          # Format\n        $replace,        # ReplaceWith\n        2                # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Text not found: $find\"\n    }\n}\n\nReplace-Text \"Systems of Water Supply and Treatment\" \"Systems of Water Treatment\"\n\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2024\" \"Ativa\u00e7\u00e3o: 01/01/2025\"\n\nReplace-Text \"Fornecer aos alunos os conhecimentos b\u00e1sicos dos sistemas de abastecimento de \u00e1gua e das tecnologias de tratamento de \u00e1gua para consumo humano.\" \"Fornecer aos alunos os conhecimentos b\u00e1sicos dos sistemas principais de tratamento de \u00e1gua para consumo humano.\"\n\nReplace-Text \"Supply the students the basic knowledge of the systems of water supply and the water treatment technologies for the human consumption.\" \"Supply the students the basic knowledge of the systems of main water treatment technologies for the human consumption.\"\n\nReplace-Text \"Sistemas de Abastecimento de \u00c1gua; Tecnologias de Tratamento de \u00c1gua;\" \"Tecnologias de Tratamento de \u00c1gua;\"\n\nReplace-Text \"Systems of Water Supply; Water Treatment Technologies;\" \"Water Treatment Technologies;\"\n\nReplace-Text \"- Tipos de dimensionamento de redes de distribui\u00e7\u00e3o de \u00e1gua;- Caracter\u00edsticas\" \"- Caracter\u00edsticas\"\n\nReplace-Text \"- Reserva\u00e7\u00e3o; - Redes de distribui\u00e7\u00e3o: tipos de rede, pe\u00e7as e \u00f3rg\u00e3os acess\u00f3rios; - Sistema de Tratamento de \u00c1gua de Ciclo Completo;\" \"- Reserva\u00e7\u00e3o; - Sistema de Tratamento de \u00c1gua de Ciclo Completo;\"\n\nReplace-Text \"- Types of dimensioning of water distribution networks;- Water characteristics\" \"- Water characteristics\"\n\nReplace-Text \"Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;\" \"Reservation;- Full Cycle Water Treatment System;\"\n"}
